# The workbook stores one weekly price record per row (rows 11-110).
# A new record was added at the top of the data (new row 11), pushing every
# existing record down by one row (old row 11 -> new row 12, ..., old row
# 110 -> new row 111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; this shifts rows 11:110 down to 12:111
# and carries the number formatting (style) of the surrounding rows along
# with it (column D keeps its date style).
$ws.Rows("11:11").Insert()

# Columns that are constant across every data row in this sheet.
$ws.Cells.Item(11, 1).Value = $ws.Cells.Item(10, 1).Value2   # Mercado ID
$ws.Cells.Item(11, 2).Value = $ws.Cells.Item(10, 2).Value2   # Mercado
$ws.Cells.Item(11, 3).Value = $ws.Cells.Item(10, 3).Value2   # Región
$ws.Cells.Item(11, 5).Value = $ws.Cells.Item(10, 5).Value2   # Codreg
$ws.Cells.Item(11, 6).Value = $ws.Cells.Item(10, 6).Value2   # Categoría ID
$ws.Cells.Item(11, 7).Value = $ws.Cells.Item(10, 7).Value2   # Categoría
$ws.Cells.Item(11, 8).Value = $ws.Cells.Item(10, 8).Value2   # Variedad
$ws.Cells.Item(11, 9).Value = $ws.Cells.Item(10, 9).Value2   # Calidad
$ws.Cells.Item(11, 18).Value = $ws.Cells.Item(10, 18).Value2 # Clasificación

# New record-specific values.
$ws.Cells.Item(11, 4).Value = 44490                             # Fecha
$ws.Cells.Item(11, 10).Value = 260                               # Volumen
$ws.Cells.Item(11, 11).Value = 5000                              # Precio mínimo
$ws.Cells.Item(11, 12).Value = 6000                              # Precio máximo
$ws.Cells.Item(11, 13).Value = 5538                              # Precio promedio ponderado
$ws.Cells.Item(11, 14).Value = "$/caja 50 unidades"              # Unidad de comercialización
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"    # Origen
$ws.Cells.Item(11, 16).Value = 111                               # Precio $/Kg
$ws.Cells.Item(11, 17).Value = 50                                # Kg o Unidades
